$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 80.25
$ws.Cells.Item(9, 9).Value = 80.25
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 80.25
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = 88.75
$ws.Cells.Item(9, 14).ClearContents()

$ws.Cells.Item(43, 8).Value = 545
$ws.Cells.Item(43, 9).Value = 300
$ws.Cells.Item(43, 10).Value = 790
$ws.Cells.Item(43, 11).Value = 300
$ws.Cells.Item(43, 12).Value = 790
$ws.Cells.Item(43, 13).Value = -231
$ws.Cells.Item(43, 14).Value = -928

$ws.Cells.Item(58, 8).Value = 63
$ws.Cells.Item(58, 9).Value = 63
$ws.Cells.Item(58, 11).Value = 189
$ws.Cells.Item(58, 13).Value = -39

$ws.Cells.Item(69, 8).Value = 4999.5
$ws.Cells.Item(69, 9).Value = 4999
$ws.Cells.Item(69, 10).Value = 5000
$ws.Cells.Item(69, 11).Value = 14997
$ws.Cells.Item(69, 12).Value = 15000
$ws.Cells.Item(69, 13).Value = -14123
$ws.Cells.Item(69, 14).Value = -16748

$ws.Cells.Item(70, 8).Value = 3274.875
$ws.Cells.Item(70, 9).Value = 3824.75
$ws.Cells.Item(70, 10).Value = 2725
$ws.Cells.Item(70, 11).Value = 11474.25
$ws.Cells.Item(70, 12).Value = 8175
$ws.Cells.Item(70, 13).Value = -11204.25
$ws.Cells.Item(70, 14).Value = -8715

$ws.Cells.Item(72, 8).Value = 4999.5
$ws.Cells.Item(72, 9).Value = 4999
$ws.Cells.Item(72, 10).Value = 5000
$ws.Cells.Item(72, 11).Value = 44991
$ws.Cells.Item(72, 12).Value = 45000
$ws.Cells.Item(72, 13).Value = -40623
$ws.Cells.Item(72, 14).Value = -53736

$ws.Cells.Item(73, 8).Value = 3274.875
$ws.Cells.Item(73, 9).Value = 3824.75
$ws.Cells.Item(73, 10).Value = 2725
$ws.Cells.Item(73, 11).Value = 11474.25
$ws.Cells.Item(73, 12).Value = 8175
$ws.Cells.Item(73, 13).Value = -10538.25
$ws.Cells.Item(73, 14).Value = -10047

$ws.Cells.Item(135, 8).Value = 683
$ws.Cells.Item(135, 9).Value = 683
$ws.Cells.Item(135, 11).Value = 6147
$ws.Cells.Item(135, 13).Value = -3612

$ws.Cells.Item(137, 8).Value = 7095
$ws.Cells.Item(137, 9).Value = 7095
$ws.Cells.Item(137, 11).Value = 21285
$ws.Cells.Item(137, 13).Value = -18735

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(38, 8).Value = 2100
$ws.Cells.Item(38, 9).Value = 2100
$ws.Cells.Item(38, 11).Value = 2100
$ws.Cells.Item(38, 13).Value = -1633

$ws.Cells.Item(132, 8).Value = 555
$ws.Cells.Item(132, 9).Value = 555
$ws.Cells.Item(132, 11).Value = 1665
$ws.Cells.Item(132, 13).Value = 865

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2431
$ws.Cells.Item(20, 9).Value = 2396.5
$ws.Cells.Item(20, 11).Value = 2396.5
$ws.Cells.Item(20, 13).Value = -2149.5

$ws.Cells.Item(134, 8).Value = 2081.5
$ws.Cells.Item(134, 9).Value = 2081.5
$ws.Cells.Item(134, 11).Value = 6244.5
$ws.Cells.Item(134, 13).Value = -3709.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 4950
$ws.Cells.Item(3, 10).Value = 4950
$ws.Cells.Item(3, 12).Value = 4950
$ws.Cells.Item(3, 14).Value = -5176

$ws.Cells.Item(25, 8).Value = 2956.3333
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 2956.3333
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 2956.3333
$ws.Cells.Item(25, 14).Value = -3304.3333
$ws.Cells.Item(25, 13).ClearContents()

$ws.Cells.Item(37, 8).Value = 5000000
$ws.Cells.Item(37, 9).Value = 5000000
$ws.Cells.Item(37, 11).Value = 5000000
$ws.Cells.Item(37, 13).Value = -4999893

$ws.Cells.Item(58, 8).Value = 1884.3334
$ws.Cells.Item(58, 9).Value = 1991.8
$ws.Cells.Item(58, 10).Value = 1347
$ws.Cells.Item(58, 11).Value = 1991.8
$ws.Cells.Item(58, 12).Value = 1347
$ws.Cells.Item(58, 13).Value = -1788.8
$ws.Cells.Item(58, 14).Value = -1753

$ws.Cells.Item(86, 8).Value = 10829.333
$ws.Cells.Item(86, 9).Value = 9994.5
$ws.Cells.Item(86, 11).Value = 9994.5
$ws.Cells.Item(86, 13).Value = -8871.5

$ws.Cells.Item(89, 8).Value = 10829.333
$ws.Cells.Item(89, 9).Value = 9994.5
$ws.Cells.Item(89, 11).Value = 49972.5
$ws.Cells.Item(89, 13).Value = -44356.5

$ws.Cells.Item(105, 8).Value = 2499
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 13).ClearContents()

$ws.Cells.Item(107, 8).Value = 828.4286
$ws.Cells.Item(107, 9).Value = 833.2222
$ws.Cells.Item(107, 10).Value = 799.6667
$ws.Cells.Item(107, 11).Value = 833.2222
$ws.Cells.Item(107, 12).Value = 799.6667
$ws.Cells.Item(107, 13).Value = 1086.7778
$ws.Cells.Item(107, 14).Value = -4639.6667

$ws.Cells.Item(136, 8).Value = 1884.3334
$ws.Cells.Item(136, 9).Value = 1991.8
$ws.Cells.Item(136, 10).Value = 1347
$ws.Cells.Item(136, 11).Value = 5975.4
$ws.Cells.Item(136, 12).Value = 4041
$ws.Cells.Item(136, 13).Value = -3425.4
$ws.Cells.Item(136, 14).Value = -9141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 274.5
$ws.Cells.Item(26, 9).Value = 250
$ws.Cells.Item(26, 10).Value = 299
$ws.Cells.Item(26, 11).Value = 750
$ws.Cells.Item(26, 12).Value = 897
$ws.Cells.Item(26, 13).Value = -462
$ws.Cells.Item(26, 14).Value = -1473

$ws.Cells.Item(121, 8).Value = 914.1667
$ws.Cells.Item(121, 10).Value = 948
$ws.Cells.Item(121, 12).Value = 2844
$ws.Cells.Item(121, 14).Value = -5464

$ws.Cells.Item(140, 8).Value = 207.5
$ws.Cells.Item(140, 9).Value = 207.5
$ws.Cells.Item(140, 11).Value = 622.5
$ws.Cells.Item(140, 13).Value = 4557.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 3000
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).ClearContents()

$ws.Cells.Item(23, 8).Value = 12
$ws.Cells.Item(23, 9).Value = 12
$ws.Cells.Item(23, 11).Value = 12
$ws.Cells.Item(23, 13).Value = 211

$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 14).ClearContents()

$ws.Cells.Item(102, 8).Value = 5226.375
$ws.Cells.Item(102, 9).Value = 4802
$ws.Cells.Item(102, 10).Value = 6499.5
$ws.Cells.Item(102, 11).Value = 4802
$ws.Cells.Item(102, 12).Value = 6499.5
$ws.Cells.Item(102, 13).Value = -3180
$ws.Cells.Item(102, 14).Value = -9743.5

$ws.Cells.Item(132, 8).Value = 5000.6665
$ws.Cells.Item(132, 9).Value = 5000.6665
$ws.Cells.Item(132, 11).Value = 15001.9995
$ws.Cells.Item(132, 13).Value = -12471.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3114.1428
$ws.Cells.Item(22, 10).Value = 4699.75
$ws.Cells.Item(22, 12).Value = 4699.75
$ws.Cells.Item(22, 14).Value = -5289.75

$ws.Cells.Item(27, 8).Value = 3114.1428
$ws.Cells.Item(27, 10).Value = 4699.75
$ws.Cells.Item(27, 12).Value = 4699.75
$ws.Cells.Item(27, 14).Value = -4913.75

$ws.Cells.Item(68, 8).Value = 3177.889
$ws.Cells.Item(68, 9).Value = 2971.7144
$ws.Cells.Item(68, 11).Value = 2971.7144
$ws.Cells.Item(68, 13).Value = -2222.7144

$ws.Cells.Item(71, 8).Value = 3177.889
$ws.Cells.Item(71, 9).Value = 2971.7144
$ws.Cells.Item(71, 11).Value = 14858.572
$ws.Cells.Item(71, 13).Value = -11114.572

$ws.Cells.Item(136, 8).Value = 16334.333
$ws.Cells.Item(136, 9).Value = 12001.5
$ws.Cells.Item(136, 11).Value = 36004.5
$ws.Cells.Item(136, 13).Value = -33454.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 5000
$ws.Cells.Item(15, 10).Value = 5000
$ws.Cells.Item(15, 12).Value = 5000
$ws.Cells.Item(15, 14).Value = -5576

$ws.Cells.Item(107, 8).Value = 525
$ws.Cells.Item(107, 9).Value = 525
$ws.Cells.Item(107, 11).Value = 1575
$ws.Cells.Item(107, 13).Value = 345

$ws.Cells.Item(132, 8).Value = 3220.3333
$ws.Cells.Item(132, 9).Value = 3220.3333
$ws.Cells.Item(132, 11).Value = 9660.999899999999
$ws.Cells.Item(132, 13).Value = -7130.999899999999
